$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.108.90"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "2.940.77"
$ws.Range("E3").Value = "  +4.46%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'353.49"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "'112.23"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.626"
$ws.Range("D10").Value = "'39.51"
$ws.Range("E11").Value = "  +3.35%  "
$ws.Range("D12").Value = "'0.136"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "'20.12"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.404.88"
$ws.Range("E14").Value = "  +4.56%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'7.76"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "2.938.51"
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("D17").Value = "'0.982"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "52.165.47"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "'7.65"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("D21").Value = "'14.22"
$ws.Range("E21").Value = "  +5.60%  "
$ws.Range("D22").Value = "0.0₃0980"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "'71.22"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("D24").Value = "'268.35"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "'2.78"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "'0.180"
$ws.Range("E26").Value = "  +11.38%  "
$ws.Range("D27").Value = "'27.02"
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "'7.08"
$ws.Range("E29").Value = "  +13.98%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.104"
$ws.Range("E30").Value = "  +15.89%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "'10.63"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "'36.98"
$ws.Range("E33").Value = "  -4.83%  "
$ws.Range("D34").Value = "'6.09"
$ws.Range("E34").Value = "  +5.37%  "
$ws.Range("D35").Value = "'53.08"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").Value = "'0.0453"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  +5.92%  "
$ws.Range("D39").Value = "'18.67"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("D42").Value = "'0.118"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").Value = "'23.43"
$ws.Range("E43").Value = "  +6.43%  "
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'2.53"
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.196.86"
$ws.Range("E46").Value = "  +2.30%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.51"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'112.62"
$ws.Range("E48").Value = "  -7.42%  "
$ws.Range("D49").Value = "'0.248"
$ws.Range("E49").Value = "  +9.93%  "
$ws.Range("D50").Value = "'0.0349"
$ws.Range("E50").Value = "  +8.05%  "
$ws.Range("D51").Value = "'0.950"
$ws.Range("E51").Value = "  -3.83%  "
